$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = -21.8692
$ws.Range("C3").Value = -11.7852
$ws.Range("E19").Value = 16.26919999999999
$ws.Range("A21").Value = -20.38019999999998
$ws.Range("A23").Value = -20.37769999999998
$ws.Range("C24").Value = -12.691
$ws.Range("E24").Value = 16.40829999999999
$ws.Range("A25").Value = -21.70029999999999
$ws.Range("B27").Value = 6.118500000000004
$ws.Range("E30").Value = 15.72109999999999
$ws.Range("B31").Value = 4.9118
$ws.Range("E31").Value = 16.36609999999999
$ws.Range("E33").Value = 17.14690000000002
$ws.Range("B39").Value = 9.613200000000003
$ws.Range("B48").Value = 5.027500000000003
$ws.Range("B51").Value = 5.231600000000001
$ws.Range("B52").Value = 5.131600000000001
$ws.Range("A53").Value = -21.75970000000001
$ws.Range("B55").Value = 6.022099999999998
$ws.Range("E55").Value = 16.57609999999999
$ws.Range("B56").Value = 4.8953
$ws.Range("A57").Value = -21.90300000000001
$ws.Range("B57").Value = 5.859700000000005
$ws.Range("C57").Value = -12.37319999999999
$ws.Range("A59").Value = -22.29350000000001
$ws.Range("C61").Value = -13.34869999999999
$ws.Range("E65").Value = 17.15270000000002
$ws.Range("A69").Value = -21.59459999999999
$ws.Range("C70").Value = -11.6604
$ws.Range("E70").Value = 17.29890000000002
$ws.Range("B73").Value = 9.040100000000004
$ws.Range("E75").Value = 16.4776
$ws.Range("A79").Value = -20.65680000000002
$ws.Range("A83").Value = -21.98129999999999
$ws.Range("E83").Value = 16.445
$ws.Range("C86").Value = -13.3951
$ws.Range("B89").Value = 5.171799999999994
$ws.Range("B90").Value = 5.561100000000001
$ws.Range("A93").Value = -21.35319999999999
$ws.Range("E96").Value = 15.9165
$ws.Range("E97").Value = 16.6048
$ws.Range("C98").Value = -11.4996
$ws.Range("C100").Value = -12.3223
$ws.Range("C102").Value = -13.1157
